# Insert a new daily price record as row 37, pushing the existing row 37
# (and everything below it) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(37).Insert()

$ws.Range("A37").Value = 2
$ws.Range("B37").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 44987
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 100112030
$ws.Range("G37").Value = "Poroto granado"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 400
$ws.Range("K37").Value = 21000
$ws.Range("L37").Value = 23000
$ws.Range("M37").Value = 22000
$ws.Range("N37").Value = '$/malla 25 kilos'
$ws.Range("O37").Value = "Provincia de Limarí"
$ws.Range("P37").Value = 880
$ws.Range("Q37").Value = 25
$ws.Range("R37").Value = "Hortaliza"
